$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates as described by the authoritative diff.
# Columns D (Price) and E (Volume/1h) contain numeric-looking text that must
# stay plain text (inline string) -- a leading apostrophe forces Excel to
# store the literal text instead of converting it to a Number/Percentage.

$ws.Range("D2").Value = "'300.93"
$ws.Range("E2").Value = "'-0.16%"
$ws.Range("E3").Value = "'2.01%"
$ws.Range("D4").Value = "'4.960"
$ws.Range("E4").Value = "'-2.59%"
$ws.Range("D5").Value = "'0.07615"
$ws.Range("E5").Value = "'-2.32%"
$ws.Range("D6").Value = "'1.947"
$ws.Range("E6").Value = "'-12.74%"
$ws.Range("D7").Value = "'7.829"
$ws.Range("E7").Value = "'0.45%"
$ws.Range("D8").Value = "'3.788"
$ws.Range("E8").Value = "'-0.74%"
$ws.Range("D9").Value = "'0.9169"
$ws.Range("E9").Value = "'-0.15%"
$ws.Range("D10").Value = "'0.1765"
$ws.Range("E10").Value = "'0.64%"
$ws.Range("D11").Value = "'0.07813"
$ws.Range("E11").Value = "'3.59%"
$ws.Range("D12").Value = "'0.08527"
$ws.Range("E12").Value = "'-5.19%"
$ws.Range("D13").Value = "'0.03165"
$ws.Range("E13").Value = "'4.57%"
$ws.Range("D14").Value = "'0.1000"
$ws.Range("E14").Value = "'-0.19%"
$ws.Range("D15").Value = "'0.001515"
$ws.Range("E15").Value = "'0.21%"
$ws.Range("D16").Value = "'0.005858"
$ws.Range("E16").Value = "'-3.24%"
$ws.Range("D18").Value = "'3.461"
$ws.Range("D19").Value = "'2.153"
$ws.Range("E19").Value = "'-4.38%"
$ws.Range("D20").Value = "'0.3346"
$ws.Range("E20").Value = "'1.63%"
$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").Value = "'-2.79%"
$ws.Range("D22").Value = "'4.267"
$ws.Range("E22").Value = "'-1.56%"
$ws.Range("E23").Value = "'9.58%"
$ws.Range("D24").Value = "'0.04513"
$ws.Range("E24").Value = "'-1.77%"
$ws.Range("E25").Value = "'-2.41%"
$ws.Range("D26").Value = "'0.004389"
$ws.Range("E26").Value = "'-1.86%"
$ws.Range("E27").Value = "'0.19%"
$ws.Range("D39").Value = "'0.01707"
$ws.Range("E39").Value = "'-3.46%"
$ws.Range("D40").Value = "'0.04677"
$ws.Range("E40").Value = "'-2.20%"
$ws.Range("D41").Value = "'0.007462"
$ws.Range("E41").Value = "'-0.87%"
$ws.Range("D42").Value = "'0.1350"
$ws.Range("E42").Value = "'-0.62%"
$ws.Range("D43").Value = "'0.002332"
$ws.Range("E43").Value = "'6.60%"
$ws.Range("D44").Value = "'0.01045"
$ws.Range("E44").Value = "'1.83%"
$ws.Range("D45").Value = "'0.00006261"
$ws.Range("E45").Value = "'0.82%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.20%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "'0.8233"
$ws.Range("E47").Value = "'12.24%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.003105"
$ws.Range("E48").Value = "'-61.13%"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.20%"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.20%"
